$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Status" column (B) for the tasks whose status progressed since
# the last save of the Sprint 2 Backlog.
$ws.Range("B2").Value  = "In Progress"   # Sub-diagram 1: To do -> In Progress
$ws.Range("B8").Value  = "Done"          # Metrics 2: In Progress -> Done
$ws.Range("B9").Value  = "Done"          # Metrics 3: In Progress -> Done
$ws.Range("B10").Value = "In Progress"   # Metrics 4: To do -> In Progress
$ws.Range("B12").Value = "Reviewing"     # Design Pattern 1: In Progress -> Reviewing
$ws.Range("B13").Value = "Reviewing"     # Design Pattern 2: In Progress -> Reviewing
$ws.Range("B14").Value = "Reviewing"     # Design Pattern 3: In Progress -> Reviewing
$ws.Range("B15").Value = "Done"          # Design Pattern 4: Reviewing -> Done
$ws.Range("B20").Value = "Done"          # Design Pattern 9: Reviewing -> Done

# Move the active selection to B8, matching where the author left off editing.
$ws.Range("B8").Select()
